# Update view-count style values on "展览" (sheet 1) and "全部类型" (sheet 4)
# Each row on "全部类型" mirrors a row on "展览" (and the other per-category
# sheets), so the same events need the same value bump applied in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 5226
$wsExhibit.Range("F19").Value = 164
$wsExhibit.Range("F22").Value = 5982
$wsExhibit.Range("F24").Value = 41
$wsExhibit.Range("F26").Value = 6296

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 5226
$wsAll.Range("F23").Value = 164
$wsAll.Range("F26").Value = 5982
$wsAll.Range("F28").Value = 41
$wsAll.Range("F30").Value = 6296
